$d = $word.ActiveDocument
$lq = [char]0x2018
$rq = [char]0x2019

# --- Paragraph 3: "Bash sophon_665k_ sq30.sh" ->
#     "Download 'Sophon-7b-pretrain-cluster' from GoogleDrive, and put it under './checkpoints/projector'"
$p3 = $d.Paragraphs.Item(3)
$newP3 = "Download " + $lq + "Sophon-7b-pretrain-cluster" + $rq + " from GoogleDrive, and put it under " + $lq + "./checkpoints/projector" + $rq
$p3.Range.Find.Execute("Bash sophon_665k_ sq30.sh", $true, $false, $false, $false, $false, $true, 1, $false, $newP3, 2) | Out-Null

# --- Paragraph 4: "Bash sophon_665k_ sq50.sh" -> "Bash sophon_665k_clu_sq.sh"
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Find.Execute("sophon_665k_ sq50", $true, $false, $false, $false, $false, $true, 1, $false, "sophon_665k_clu_sq", 2) | Out-Null

# --- Paragraph 5: "Bash sophon_665k_ sq50_weightdecay.sh" -> "Bash pretrain.sh"
#     then two new list paragraphs appended after it.
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Find.Execute("sophon_665k_ sq50_weightdecay", $true, $false, $false, $false, $false, $true, 1, $false, "pretrain", 2) | Out-Null

$newPara1 = "After finishing step 3, put " + $lq + "./checkpoints/Sophon-7b-pretrain-qav" + $rq + " under " + $lq + "./checkpoints/projector" + $rq
$newPara2 = "Bash sophon_665k_ sq30.sh"

$p5.Range.InsertParagraphAfter() | Out-Null
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Text = $newPara1

$p6.Range.InsertParagraphAfter() | Out-Null
$p7 = $d.Paragraphs.Item(7)
$p7.Range.Text = $newPara2
